# Apply updated cryptocurrency price/volume figures (Price & Volume(1h) columns)
# to the worksheet, matching the latest scrape performed by the GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.115.74"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "1.636.39"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'216.85"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.253"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").Value = "'19.91"
$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "1.865.45"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "1.636.41"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").Value = "'0.540"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "'66.67"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "27.116.95"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "'216.73"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "'6.84"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").Value = "'2.51"
$ws.Range("E23").Value = "  +3.01%  "
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").Value = "'146.56"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'7.39"
$ws.Range("E27").Value = "  +2.07%  "
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").Value = "'15.66"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("D34").Value = "1.302.47"
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("D38").Value = "'0.854"
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("E39").Value = "  +1.74%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'0.806"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "'2.24"
$ws.Range("E42").Value = "  +6.31%  "
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("D44").Value = "1.775.96"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "'61.72"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").Value = "'91.23"
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "'7.61"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("E51").Value = "  -0.37%  "
